# Generate Report for Handback
# Updates the zh-cn and de-de localization-status sheets:
#  - Status column: "Ready for handoff" -> "Handed back: in sync with en-US"
#  - Adds "Latest Target File" (F) and "Latest Handback File" (G) hyperlinked values
#  - Sets the real "Latest Handback DateTime" (H) instead of the 0001-01-01 placeholder

$wb = $excel.ActiveWorkbook

$mdUrl1  = "https://github.com/OpenLocalizationTest/oltest/blob/059e5c38ee30241430aafe24597291d0732bac4c/e2e/73c39128-3f88-4cbe-9850-3efc160b9f93.md"
$mdName1 = "73c39128-3f88-4cbe-9850-3efc160b9f93.md"

function Update-LocSheet {
    param($SheetName, $XlfUrl, $XlfName, $HandbackDateTime)

    $ws = $wb.Worksheets.Item($SheetName)

    # Status -> Handed back
    $ws.Range("C2").Value = "Handed back: in sync with en-US"
    $ws.Range("C3").Value = "Handed back: in sync with en-US"

    # Latest Target File (F) + Latest Handback File (G), rows 2 and 3
    foreach ($row in 2, 3) {
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 6), $mdUrl1, "", "", $mdName1)
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 7), $XlfUrl, "", "", $XlfName)
    }

    # Latest Handback DateTime (H)
    $ws.Range("H2").Value = $HandbackDateTime
    $ws.Range("H3").Value = $HandbackDateTime
}

$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7e7ffafbc2c2ec44a8fac140516db6d92d3f2e42/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/73c39128-3f88-4cbe-9850-3efc160b9f93.9938c5df5a8f32b29a7ede5bc650f7b859d603c6.zh-cn.xlf"
$zhXlfName = "73c39128-3f88-4cbe-9850-3efc160b9f93.9938c5df5a8f32b29a7ede5bc650f7b859d603c6.zh-cn.xlf"

$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d0d6509bdd07110534558b36c9446faef4273ab8/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/73c39128-3f88-4cbe-9850-3efc160b9f93.9938c5df5a8f32b29a7ede5bc650f7b859d603c6.de-de.xlf"
$deXlfName = "73c39128-3f88-4cbe-9850-3efc160b9f93.9938c5df5a8f32b29a7ede5bc650f7b859d603c6.de-de.xlf"

Update-LocSheet "zh-cn" $zhXlfUrl $zhXlfName "2016-03-20 16:57:59"
Update-LocSheet "de-de" $deXlfUrl $deXlfName "2016-03-20 16:58:06"

"done"
